# Select/activate "Sheet2" (the first tab) and stamp a header label into A1.
# This mirrors the author's change: the workbook's active tab moves from
# "Sheet1" (tab index 1) to "Sheet2" (tab index 0), and a new label
# "Break Out" is written into Sheet2!A1 (adds a shared string + extends the
# sheet's used range to A1:A4).

$wb = $excel.ActiveWorkbook

$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Activate()

$ws2.Range("A1").Value = "Break Out"
